$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = 133297
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Archivo dañado"

$ws.Range("A11").Value = 98617
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "098619 es la misma canción y mismo fragmento"

$ws.Range("A12").Value = 98565
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Dura 1 segundo"

$ws.Range("A13").Value = 98567
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Dura 1 segundo"

$ws.Range("A14").Value = 98569
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Dura 1 segundo"

$ws.Range("A15").Value = 72059
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "No tiene sonido alguno"

$ws.Range("C15").Select()
